# "dodanie punktów odniesienia i wag" - adding a reference point / weight
# on the Arkusz1 route-scoring sheet: avalanche-risk field (B11) goes
# from 0 to 1, which ripples into the computed time estimate (D11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

$ws.Range("B11").Value = 1

# Leave the selection where the author ended up after the edit.
$ws.Activate()
$ws.Range("C14").Select()
